$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# The export previously included a placeholder row for 2025-08-25 with no
# data. Remove that row; everything below shifts up to close the gap.
$ws.Rows.Item(2).Delete()
